$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the Price (D) column cells we are updating as Text first, so that
# purely-numeric-looking values (e.g. "69.00", "7.10") are preserved exactly
# as literal text, matching the source data (inline/shared strings), instead
# of being auto-converted into numbers by Excel (which would drop the exact
# formatting, e.g. "69.00" -> 69).
$dCells = @("D2", "D3", "D5", "D6", "D14", "D15", "D16", "D17", "D18", "D19", "D22", "D23", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D36", "D37", "D38", "D42", "D43", "D45", "D49", "D50", "D51")
foreach ($cellRef in $dCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "45.432.60"
$ws.Range("E2").Value = "  +2.85%  "
$ws.Range("D3").Value = "2.426.70"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "318.64"
$ws.Range("E5").Value = "  +3.55%  "
$ws.Range("D6").Value = "103.05"
$ws.Range("E6").Value = "  +2.87%  "
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("E9").Value = "  +6.13%  "
$ws.Range("E10").Value = "  +0.76%  "
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("E13").Value = "  -2.80%  "
$ws.Range("D14").Value = "7.10"
$ws.Range("E14").Value = "  +2.36%  "
$ws.Range("D15").Value = "2.806.65"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("D16").Value = "2.417.76"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").Value = "0.843"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "45.342.75"
$ws.Range("E18").Value = "  +2.65%  "
$ws.Range("D19").Value = "12.23"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("E21").Value = "  +2.00%  "
$ws.Range("D22").Value = "69.00"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").Value = "244.67"
$ws.Range("E24").Value = "  -1.42%  "
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "25.76"
$ws.Range("E27").Value = "  +1.81%  "
$ws.Range("D28").Value = "2.28"
$ws.Range("E28").Value = "  -2.72%  "
$ws.Range("D29").Value = "9.61"
$ws.Range("E29").Value = "  +1.21%  "
$ws.Range("D30").Value = "49.43"
$ws.Range("E30").Value = "  +2.62%  "
$ws.Range("D31").Value = "32.95"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "20.30"
$ws.Range("E32").Value = "  +8.89%  "
$ws.Range("D33").Value = "0.125"
$ws.Range("E33").Value = "  +5.18%  "
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").Value = "0.0766"
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "4.47"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "1.87"
$ws.Range("E38").Value = "  -2.63%  "
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("E40").Value = "  -2.44%  "
$ws.Range("E41").Value = "  +0.69%  "
$ws.Range("D42").Value = "2.20"
$ws.Range("E42").Value = "  -3.68%  "
$ws.Range("D43").Value = "20.62"
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("E44").Value = "  +0.88%  "
$ws.Range("D45").Value = "1.939.31"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("E46").Value = "  -2.95%  "
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("E48").Value = "  +9.05%  "
$ws.Range("D49").Value = "9.14"
$ws.Range("E49").Value = "  -2.56%  "
$ws.Range("D50").Value = "76.98"
$ws.Range("E50").Value = "  +4.46%  "
$ws.Range("D51").Value = "4.79"
$ws.Range("E51").Value = "  +5.73%  "
